$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = "28.341.42"
$ws.Cells.Item(2, 5).Value = "  +3.75%  "

# Row 3
$ws.Cells.Item(3, 4).Value = "1.820.19"
$ws.Cells.Item(3, 5).Value = "  +4.61%  "

# Row 4
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = "1.001"
$ws.Cells.Item(4, 5).Value = "  -0.31%  "

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "329.09"
$ws.Cells.Item(5, 5).Value = "  +2.53%  "

# Row 6
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "0.9992"
$ws.Cells.Item(6, 5).Value = "  -0.22%  "

# Row 7
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "0.4377"
$ws.Cells.Item(7, 5).Value = "  +4.94%  "

# Row 8
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "0.3697"
$ws.Cells.Item(8, 5).Value = "  +3.56%  "

# Row 9
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "45.06"
$ws.Cells.Item(9, 5).Value = "  -0.38%  "

# Row 10
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "0.07729"
$ws.Cells.Item(10, 5).Value = "  +4.93%  "

# Row 11
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "1.144"
$ws.Cells.Item(11, 5).Value = "  +3.33%  "

# Row 12
$ws.Cells.Item(12, 2).Value = "Solana"
$ws.Cells.Item(12, 3).Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "22.29"
$ws.Cells.Item(12, 5).Value = "  +4.58%  "

# Row 13
$ws.Cells.Item(13, 2).Value = "BinanceUSD"
$ws.Cells.Item(13, 3).Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "0.9991"
$ws.Cells.Item(13, 5).Value = "  -0.34%  "

# Row 14
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "6.334"
$ws.Cells.Item(14, 5).Value = "  +4.82%  "

# Row 15
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "7.618"
$ws.Cells.Item(15, 5).Value = "  +6.56%  "

# Row 16
$ws.Cells.Item(16, 4).Value = "1.832.83"
$ws.Cells.Item(16, 5).Value = "  +5.98%  "

# Row 17
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "93.74"
$ws.Cells.Item(17, 5).Value = "  +7.86%  "

# Row 18
$ws.Cells.Item(18, 5).Value = "  +2.17%  "

# Row 19
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "0.06523"
$ws.Cells.Item(19, 5).Value = "  +8.36%  "

# Row 20
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "0.9994"
$ws.Cells.Item(20, 5).Value = "  -0.25%  "

# Row 21
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "17.55"
$ws.Cells.Item(21, 5).Value = "  +4.74%  "

# Row 22
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "6.281"
$ws.Cells.Item(22, 5).Value = "  +3.66%  "

# Row 23
$ws.Cells.Item(23, 4).Value = "28.373.21"
$ws.Cells.Item(23, 5).Value = "  +3.71%  "

# Row 24
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "11.71"
$ws.Cells.Item(24, 5).Value = "  +3.49%  "

# Row 25
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "2.129"
$ws.Cells.Item(25, 5).Value = "  -9.21%  "

# Row 26
$ws.Cells.Item(26, 2).Value = "Monero"
$ws.Cells.Item(26, 3).Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "160.19"
$ws.Cells.Item(26, 5).Value = "  +4.97%  "

# Row 27
$ws.Cells.Item(27, 2).Value = "EthereumClassic"
$ws.Cells.Item(27, 3).Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "20.86"
$ws.Cells.Item(27, 5).Value = "  +3.02%  "

# Row 28
$ws.Cells.Item(28, 4).Value = "2.034.08"
$ws.Cells.Item(28, 5).Value = "  +4.81%  "

# Row 29
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "2.321"
$ws.Cells.Item(29, 5).Value = "  -1.96%  "

# Row 30
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "129.45"
$ws.Cells.Item(30, 5).Value = "  +3.20%  "

# Row 31
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "1.219"
$ws.Cells.Item(31, 5).Value = "  +4.54%  "

# Row 32
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "5.965"
$ws.Cells.Item(32, 5).Value = "  +5.49%  "

# Row 33
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "0.09225"
$ws.Cells.Item(33, 5).Value = "  +1.67%  "

# Row 34
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "3.615"
$ws.Cells.Item(34, 5).Value = "  +0.44%  "

# Row 35
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "12.96"
$ws.Cells.Item(35, 5).Value = "  +2.52%  "

# Row 36
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "0.02370"
$ws.Cells.Item(36, 5).Value = "  +4.72%  "

# Row 37
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "0.2193"
$ws.Cells.Item(37, 5).Value = "  +3.19%  "

# Row 38
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "5.203"
$ws.Cells.Item(38, 5).Value = "  +3.25%  "

# Row 39
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "0.6611"
$ws.Cells.Item(39, 5).Value = "  +4.28%  "

# Row 40
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "0.06227"
$ws.Cells.Item(40, 5).Value = "  +3.41%  "

# Row 41
$ws.Cells.Item(41, 2).Value = "FraxShare"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "8.221"
$ws.Cells.Item(41, 5).Value = "  +3.91%  "

# Row 42
$ws.Cells.Item(42, 2).Value = "TrustWalletToken"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "1.196"
$ws.Cells.Item(42, 5).Value = "  +0.96%  "

# Row 43
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "1.426"
$ws.Cells.Item(43, 5).Value = "  -0.07%  "

# Row 44
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "0.9991"
$ws.Cells.Item(44, 5).Value = "  -0.24%  "

# Row 45
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "13.93"
$ws.Cells.Item(45, 5).Value = "  +2.19%  "

# Row 46
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "0.6164"
$ws.Cells.Item(46, 5).Value = "  +6.86%  "

# Row 47
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "3.761"
$ws.Cells.Item(47, 5).Value = "  +1.77%  "

# Row 48
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "126.84"
$ws.Cells.Item(48, 5).Value = "  +1.48%  "

# Row 49
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "2.036"
$ws.Cells.Item(49, 5).Value = "  +5.37%  "

# Row 50
$ws.Cells.Item(50, 2).Value = "Cronos"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "0.07031"
$ws.Cells.Item(50, 5).Value = "  +3.28%  "

# Row 51
$ws.Cells.Item(51, 2).Value = "EOS"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "1.161"
$ws.Cells.Item(51, 5).Value = "  +5.95%  "
